$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the "Затр. время в сек." shared formula's styled range down to row 15
# (mirrors E2:E14's bold formatting; cell itself stays empty, no D15 to divide)
$ws.Range("E15").Font.Bold = $true

# New row 17: "Между запусками в 1 заказе"
$ws.Range("A17").Value = "Между запусками в 1 заказе"
$ws.Range("E17").Value = 10
$ws.Range("D17").Formula = "=E17*B18"

# Row 18: add elapsed-time style cell + totals
$ws.Range("D18").NumberFormat = "mm:ss"
$ws.Range("E18").Formula = "=SUM(E2:E15)"
$ws.Range("F18").Formula = "=SUM(E4:E14)"

# Row 19: additional total
$ws.Range("F19").Formula = "=SUM(E5:E14)"

# Match the final selection recorded in the workbook
$ws.Range("D17").Select()
